$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # GIN
$ws2 = $wb.Worksheets.Item(2)   # BitBucket Repos

# ---------------------------------------------------------------------------
# Sheet "GIN" (sheet1): column A changes to the same repository path on every
# row, and a new row 4 for "Wichita_Tag_Printer" is appended (with hyperlink).
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "FEI_SHARED\Repository\CNH Grand Island"
$ws1.Range("A3").Value = "FEI_SHARED\Repository\CNH Grand Island"

$ws1.Range("A4").Value = "FEI_SHARED\Repository\CNH Grand Island"
$ws1.Range("B4").Value = "\\s1cn1faras14\SharedData\Temp\repo\Wichita_Tag_Printer"
$ws1.Range("C4").Value = "X"
$ws1.Range("D4").Value = "Wichita_Tag_Printer"

$ws1.Hyperlinks.Add($ws1.Range("B4"), "file:///\\s1cn1faras14\SharedData\Temp\repo\Wichita_Tag_Printer") | Out-Null

# Match the styling used by the other two hyperlink cells (B2/B3) instead of
# leaving the brand new style entry that Hyperlinks.Add() just created.
$ws1.Range("B2").Copy()
$ws1.Range("B4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column widths widened on sheet1 (target raw OOXML width = 45 / 56; the
# COM ColumnWidth property adds a ~0.8333 padding offset on save, so we
# dial in a value that round-trips to the exact widths we need).
$ws1.Columns.Item(1).ColumnWidth = 44.166666666666664
$ws1.Columns.Item(2).ColumnWidth = 55.166666666666664

# ---------------------------------------------------------------------------
# Sheet "BitBucket Repos" (sheet2): row 2 now describes Wichita_Tag_Printer,
# and two more rows (MODIS, UTS2018) are appended.
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "Wichita_Tag_Printer"
$ws2.Range("B2").Value = "D:\\Temp\\Wichita_Tag_Printer"
$ws2.Range("C2").Value = "FEI_SHARED\Repository\BitBucket Repos\Wichita_Tag_Printer"

$ws2.Range("A3").Value = "MODIS"
$ws2.Range("B3").Value = "D:\Temp\MODIS"
$ws2.Range("C3").Value = "FEI_SHARED\Repository\BitBucket Repos\MODIS"

$ws2.Range("A4").Value = "UTS2018"
$ws2.Range("B4").Value = "D:\Temp\UTS"
$ws2.Range("C4").Value = "FEI_SHARED\Repository\Bitbucket Repos\UTS"

# ---------------------------------------------------------------------------
# View state: GIN becomes the active/selected sheet with A4 selected, while
# BitBucket Repos loses the tab-selected flag but keeps a selection on A4.
# ---------------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("A4").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("A4").Select() | Out-Null
